# Updated data to reflect new requirement separation
# The "Terms Typically Offered" column (D) is pushed out to G, and three new
# columns are inserted in between: Corequisites (D), Concurrent (E) and
# Recommended (F). Row 8 had its corequisite text jammed into the old D
# value ("F, SPCorequisite: GEOL 102 or GEOL 201.") - that gets split back
# into the real Terms value ("F, SP" -> G8) and the real corequisite text
# ("GEOL 102 or GEOL 201." -> D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "Terms Typically Offered" column (D) values for the
# 23 data/header rows before we overwrite anything.
$lastRow = 23
$oldTerms = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $oldTerms[$r] = $ws.Cells.Item($r, 4).Value2
}

# Row 8's old D value mixes the terms-offered text with the corequisite
# text; split it out explicitly. The source data uses a non-breaking space
# (U+00A0) between a course prefix and its number throughout the sheet, so
# keep that convention for the text we are relocating.
$nbsp = [char]0x00A0
$row8Coreq = "GEOL" + $nbsp + "102 or GEOL" + $nbsp + "201."
$row8Terms = "F, SP"

# Header row: insert the three new column headers, and move the old header
# text for column D out to the new column G.
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"
$ws.Cells.Item(1, 7).Value = $oldTerms[1]

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 8) {
        $ws.Cells.Item($r, 4).Value = $row8Coreq
        $ws.Cells.Item($r, 5).Value = "NA"
        $ws.Cells.Item($r, 6).Value = "NA"
        $ws.Cells.Item($r, 7).Value = $row8Terms
    } else {
        $ws.Cells.Item($r, 4).Value = "NA"
        $ws.Cells.Item($r, 5).Value = "NA"
        $ws.Cells.Item($r, 6).Value = "NA"
        $ws.Cells.Item($r, 7).Value = $oldTerms[$r]
    }
}

# Row 11's prerequisite text changed ("one of the following:" -> "one of the"),
# again preserving the sheet's non-breaking-space convention.
$ws.Cells.Item(11, 3).Value = "GEOL" + $nbsp + "301 and one of the STAT" + $nbsp + "217, STAT" + $nbsp + "218, STAT" + $nbsp + "301, STAT" + $nbsp + "312, or STAT" + $nbsp + "321."

Write-Output "done"
